$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new forecast row (54), mirroring the layout/formatting of the
# preceding data row (date in col A formatted like the rows above it).
$newRow = 54
$prevRow = $newRow - 1

$ws.Cells.Item($prevRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = 45986
$ws.Cells.Item($newRow, 2).Value = 2025
$ws.Cells.Item($newRow, 3).Value = 0.8976398032236155
$ws.Cells.Item($newRow, 4).Value = 2026
$ws.Cells.Item($newRow, 5).Value = 0.4275768375374467
